# Update "想去人数" (F) and "最低票价" (G) values on the two sheets that
# contain the full event listing: "展览" (sheet 1) and "全部类型" (sheet 4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 206
    $ws.Range("F4").Value = 5
    $ws.Range("F5").Value = 6515
    $ws.Range("F6").Value = 77
    $ws.Range("F7").Value = 427
    $ws.Range("F8").Value = 131
    $ws.Range("F9").Value = 5882
    $ws.Range("G9").Value = 68
    $ws.Range("F10").Value = 37
    $ws.Range("F12").Value = 1232
    $ws.Range("F16").Value = 87
    $ws.Range("F18").Value = 334
    $ws.Range("F21").Value = 4198

    if ($sheetName -eq "展览") {
        $ws.Range("F23").Value = 182
    } else {
        $ws.Range("F24").Value = 182
    }
}
